$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the new daily entries (day 20 and day 21 of 08/2025),
# pushing the existing rows 21+ down to 23+.
$ws.Rows("21:22").Insert()

# Update existing row 20 (day 19) total_venda value.
$ws.Range("B20").Value = 16035.52

# Fill the newly inserted row 21 (day 20, 08/2025).
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 24420.32
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 2025
$ws.Range("E21").Value = "08/2025"

# Fill the newly inserted row 22 (day 21, 08/2025).
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 20734.7
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 2025
$ws.Range("E22").Value = "08/2025"
